$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B2 (Farm_ID) value from 100 to 104
$ws.Range("B2").Value = 104

# Update B11 (N_Total) value from 1.8 to 5
$ws.Range("B11").Value = 5

# Move the active selection to B11 (matches author's final cursor position)
$ws.Range("B11").Select()
